# add end switches handlers
#
# This script reproduces, via Excel COM automation, the changes described
# by the commit "add end switches handlers":
#   - InputRegs: trim 4 unused trailing blank rows (19-22), reset view
#   - HoldingRegs: re-color rows 53-62 to the normal (non-highlighted) look,
#                  trim 3 unused trailing blank rows (63-65), deselect tab
#   - DiscreteInputs: untouched
#   - Coils: insert a new "Save Cfg" register row, renumber the reserved
#            range text from "0-63" to "1-63", make this sheet the active tab

$wb = $excel.ActiveWorkbook

$wsInput    = $wb.Worksheets.Item(1)   # InputRegs
$wsHolding  = $wb.Worksheets.Item(2)   # HoldingRegs
$wsDiscrete = $wb.Worksheets.Item(3)   # DiscreteInputs
$wsCoils    = $wb.Worksheets.Item(4)   # Coils

# ---------------------------------------------------------------------
# InputRegs: remove the 4 trailing empty rows (19:22) and reset the view
# ---------------------------------------------------------------------
$wsInput.Range("A19:E22").EntireRow.Delete()

# ---------------------------------------------------------------------
# HoldingRegs: rows 53-62 lose their special highlight fill and fall back
# to the plain "normal" row look used elsewhere on the sheet (same look
# as row 52); then drop the 3 trailing empty rows (63:65)
# ---------------------------------------------------------------------
$wsHolding.Range("A52").Copy()
$wsHolding.Range("A53:F62").PasteSpecial(-4122)
$wsHolding.Range("A63:F65").EntireRow.Delete()

# ---------------------------------------------------------------------
# Coils: insert a new row 2 for the "Save Cfg" register
# ---------------------------------------------------------------------
$wsCoils.Range("A2:E2").EntireRow.Insert()

# give the new row the same highlighted look as the sheet's other
# "register id" rows (e.g. InputRegs row 2)
$wsInput.Range("A2:E2").Copy()
$wsCoils.Range("A2:E2").PasteSpecial(-4122)

$wsCoils.Range("A2").Value = 0
$wsCoils.Range("B2").Value = "Save Cfg"
$wsCoils.Range("C2").Value = ""
$wsCoils.Range("D2").Value = ""
$wsCoils.Range("E2").Value = "Сохранить настройки устройства во FLASH"

# the old reserved-range row (now pushed down to row 3) changes its
# description from "0-63" to "1-63"
$wsCoils.Range("A3").Value = "1-63"

# ---------------------------------------------------------------------
# View / selection state
# ---------------------------------------------------------------------
[void]$wsInput.Activate()
[void]$wsInput.Range("B13").Select()

[void]$wsHolding.Activate()
[void]$wsHolding.Range("D13").Select()

[void]$wsDiscrete.Activate()
[void]$wsDiscrete.Range("H7").Select()

[void]$wsCoils.Activate()
[void]$wsCoils.Range("H10").Select()
